$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for every data row (2-44) is being bumped from
# serial date 45179 (2023-09-10) to 45180 (2023-09-11).
$lastRow = 44
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45180
}
